$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update revised case counts (column C) for rows 290-318 ---
$newC = @{
    290 = 2989; 291 = 2608; 292 = 1712; 293 = 1193; 294 = 3522; 295 = 3135;
    296 = 2923; 297 = 3001; 298 = 2850; 299 = 1767; 300 = 1187; 301 = 3593;
    302 = 3794; 303 = 2940; 304 = 446;  305 = 3364; 306 = 2903; 307 = 1753;
    308 = 5484; 309 = 5842; 310 = 6051; 311 = 5740; 312 = 5285; 313 = 2166;
    314 = 2092; 315 = 6035; 316 = 4918; 317 = 4131; 318 = 3640
}

foreach ($r in $newC.Keys) {
    $ws.Cells.Item($r, 3).Value = $newC[$r]
}

# --- 2. Append two new rows (319 and 320) with the latest data ---

# Row 319 - copy formatting from row 318, then set values/formulas
$ws.Range("A318").Copy($ws.Range("A319"))
$ws.Range("A319").Value = 44176
$ws.Range("B318").Copy($ws.Range("B319"))
$ws.Range("B319").Formula = "=C319+B318"
$ws.Range("C318").Copy($ws.Range("C319"))
$ws.Range("C319").Value = 2597
$ws.Range("D318").Copy($ws.Range("D319"))
$ws.Range("D319").Formula = "=AVERAGE(C313:C319)"

# Row 320
$ws.Range("A319").Copy($ws.Range("A320"))
$ws.Range("A320").Value = 44177
$ws.Range("B319").Copy($ws.Range("B320"))
$ws.Range("B320").Formula = "=C320+B319"
$ws.Range("C319").Copy($ws.Range("C320"))
$ws.Range("C320").Value = 211
$ws.Range("D319").Copy($ws.Range("D320"))
$ws.Range("D320").Formula = "=AVERAGE(C314:C320)"

# --- 3. Update sheet view (frozen pane / selection) to match latest scroll position ---
$excel.ActiveWindow.ScrollRow = 290
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H315").Select() | Out-Null

Write-Host ("Dimension used range: " + $ws.UsedRange.Address())
Write-Host ("B320 value: " + $ws.Range("B320").Value2)
Write-Host ("D320 value: " + $ws.Range("D320").Value2)
